# İş Takip Güncellemesi - 15.12.2025 13:07:41
#
# The tracker's reference dates were all off by one day, so every
# "başlama/bitiş" style date on the two tracking sheets gets rolled back
# by exactly one day. The dates are stored as plain text (not real Excel
# dates), so each cell is forced to Text format before the new value is
# written back - this keeps them as literal "yyyy-MM-dd" strings instead
# of letting Excel re-interpret them as date serials.

$wb = $excel.ActiveWorkbook

function Shift-DateCellBack {
    param($cell)
    $s = $cell.Value2
    if ($s -ne $null -and $s -ne "") {
        $old = [DateTime]::ParseExact($s.ToString(), "yyyy-MM-dd", $null)
        $new = $old.AddDays(-1)
        $cell.NumberFormat = "@"
        $cell.Value2 = $new.ToString("yyyy-MM-dd")
    }
}

# --- Sheet "İş Takip Listesi": columns J (İŞE BAŞLAMA/YER TESLİMİ) and
#     K (İHALE BİTİŞ TARİHİ) shift back one day for rows 2-10 and 33-122.
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")
$ws1Rows = @(2..10) + @(33..122)
foreach ($r in $ws1Rows) {
    foreach ($col in @("J","K")) {
        Shift-DateCellBack $ws1.Range("$col$r")
    }
}

# --- Sheet "Güncelleme": columns I, J, N, P shift back one day for any
#     row (2-29) where the cell actually holds a date.
$ws2 = $wb.Worksheets.Item("Güncelleme")
for ($r = 2; $r -le 29; $r++) {
    foreach ($col in @("I","J","N","P")) {
        Shift-DateCellBack $ws2.Range("$col$r")
    }
}

# Row 6's "KOMİSYON DURUM" (O6) no longer reads "Yapıldı" - clear it out.
$ws2.Range("O6").Value2 = ""
